$wb = $excel.ActiveWorkbook

# --- Branches sheet: insert a new row of data (git push -u origin dev:dev1) ---
$branches = $wb.Worksheets.Item("Branches")

# Insert a new blank row above the old row 10 ("git checkout master" / ...),
# shifting it (and everything below) down by one row.
$branches.Rows.Item(10).Insert()

# Fill the freshly inserted row 9 with the new command/description pair.
$branches.Range("A9").Value = "git push -u origin dev:dev1"
$branches.Range("B9").Value = "pushes local dev branch to remote dev1 branch"

# Update the Branches sheet selection to C9 (matches the recorded selection).
$branches.Range("C9").Select()

# --- Commands sheet: no longer the active/selected tab, selection unchanged ---
# (Its selection B39 stays as-is; simply not re-activating it removes
# tabSelected from its view once another sheet is activated below.)

# --- Sheet1: becomes the active tab; keep its existing selection (B17) ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("B17").Select()
